$d = $word.ActiveDocument

# --- Edit 1: remove the stray _GoBack bookmark from the abstract paragraph (para 4) ---
$p4 = $d.Paragraphs.Item(4)
$p4xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:i/></w:rPr></w:pPr><w:r><w:rPr><w:i/></w:rPr><w:tab/><w:t xml:space="preserve">This one </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>paragragh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> description of our paper will need to be written after the rest of the paper is written, so this is just a placeholder ha </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>ha</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>ha</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t>…</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> also note this text appears different than the rest of the paper.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p4.Range.InsertXML($p4xml)

# --- Edit 2: write the two new algorithm paragraphs before "A short description..." (para 11) ---
$algTarget = $d.Paragraphs.Item(11)
$algTarget.Range.InsertParagraphBefore()
$algHolder = $d.Paragraphs.Item(11)
$algxml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t>There are two different algorithms used for these experiments.  One algorithm will take the target string and use its length to generate a population of random strings of that length.  This population is then measured against the target using the fitness function, the individuals that are selected are then mutated, with a small selection undergoing crossover.  This corresponds to asexual and sexual reproduction</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>T</w:t></w:r><w:r><w:t xml:space="preserve">he mutated individual is kept in the population as well as its mutation, while the crossover selects a random crossover point and returns the four possible combinations of the two individuals.  For example if the individuals are {ABC, </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>abc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>} and the crossover point divides them after the first character, the results would be {</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Abc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>aBC</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bcA</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>BCa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r><w:r><w:t xml:space="preserve">.  The mutation selects a random character from the string and changes it to either the character before or after it.  </w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">The other algorithm used differs in a few important ways.  The population generated consists of various length binary strings.  The binary string is then used to generate a character string.  The character string generated </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">is what the fitness test ‘grades’.  The fitness test grades the strings based on length and content.  The closer they are to the target the better their grade is.  </w:t></w:r><w:r><w:t xml:space="preserve">A selection of the population is made which will breed.  The breeding is random among the selection, and consists of a crossover similar to the other algorithm, except the crossover point is limited to the length of the shorter of the two strings as the length is not predetermined. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$algHolder.Range.InsertXML($algxml)

# --- Edit 3: strike-through the placeholder "Algorithm" section paragraphs and re-wrap them with the _GoBack bookmark ---
$s1 = $d.Paragraphs.Item(13)
$s2 = $d.Paragraphs.Item(14)
$s3 = $d.Paragraphs.Item(15)
$s1.Range.Font.StrikeThrough = 1
$s2.Range.Font.StrikeThrough = 1
$s3.Range.Font.StrikeThrough = 1
$bmRange = $d.Range($s1.Range.Start, $s3.Range.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- Edit 4: drop the lastRenderedPageBreak before "References should be mentioned somewhere in the paper" ---
$refPara = $d.Paragraphs.Item(43)
$refxml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>References should be mentioned somewhere in the paper</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$refPara.Range.InsertXML($refxml)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
